# Update Sema3a-Plxna1 LR-pair sheet with new TPM-derived metrics.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 (ECs -> ECs) ---
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.250631
$ws.Range("H2").Value = 0.751893
$ws.Range("I2").Value = 0.2648339568266264
$ws.Range("J2").Value = 0.2648339568266264
$ws.Range("M2").Value = 4.224096333333333
$ws.Range("N2").Value = 12.672289
$ws.Range("O2").Value = 0.1277189908446358
$ws.Range("P2").Value = 0.1277189908446358
$ws.Range("Q2").Value = 1.058689488119666
$ws.Range("R2").Value = 9.528205393077
$ws.Range("S2").Value = 0.03382432570728858
$ws.Range("T2").Value = 0.03382432570728858

# --- Row 3 (ECs -> FAPs) ---
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.250631
$ws.Range("H3").Value = 0.751893
$ws.Range("I3").Value = 0.2648339568266264
$ws.Range("J3").Value = 0.2648339568266264
$ws.Range("O3").Value = 0.4492078640046304
$ws.Range("P3").Value = 0.4492078640046304
$ws.Range("Q3").Value = 3.723578149634
$ws.Range("R3").Value = 33.512203346706
$ws.Range("S3").Value = 0.1189654960619833
$ws.Range("T3").Value = 0.1189654960619833

# --- Row 4 (ECs -> MuSCs) ---
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.250631
$ws.Range("H4").Value = 0.751893
$ws.Range("I4").Value = 0.2648339568266264
$ws.Range("J4").Value = 0.2648339568266264
$ws.Range("O4").Value = 0.4230731451507339
$ws.Range("P4").Value = 0.4230731451507338
$ws.Range("Q4").Value = 3.506941986581
$ws.Range("R4").Value = 31.562477879229
$ws.Range("S4").Value = 0.1120441350573545
$ws.Range("T4").Value = 0.1120441350573545

# --- Row 5 (FAPs -> ECs) ---
$ws.Range("I5").Value = 0.2480790641859371
$ws.Range("J5").Value = 0.2480790641859371
$ws.Range("M5").Value = 4.224096333333333
$ws.Range("N5").Value = 12.672289
$ws.Range("O5").Value = 0.1277189908446358
$ws.Range("P5").Value = 0.1277189908446358
$ws.Range("Q5").Value = 0.9917108086262222
$ws.Range("R5").Value = 8.925397277636
$ws.Range("S5").Value = 0.03168440772750953
$ws.Range("T5").Value = 0.03168440772750953

# --- Row 6 (FAPs -> FAPs) ---
$ws.Range("I6").Value = 0.2480790641859371
$ws.Range("J6").Value = 0.2480790641859371
$ws.Range("O6").Value = 0.4492078640046304
$ws.Range("P6").Value = 0.4492078640046304
$ws.Range("S6").Value = 0.1114390665272324
$ws.Range("T6").Value = 0.1114390665272324

# --- Row 7 (FAPs -> MuSCs) ---
$ws.Range("I7").Value = 0.2480790641859371
$ws.Range("J7").Value = 0.2480790641859371
$ws.Range("O7").Value = 0.4230731451507339
$ws.Range("P7").Value = 0.4230731451507338
$ws.Range("S7").Value = 0.1049555899311952
$ws.Range("T7").Value = 0.1049555899311952

# --- Row 8 (MuSCs -> ECs) ---
$ws.Range("G8").Value = 0.4609646666666666
$ws.Range("I8").Value = 0.4870869789874365
$ws.Range("J8").Value = 0.4870869789874365
$ws.Range("M8").Value = 4.224096333333333
$ws.Range("N8").Value = 12.672289
$ws.Range("O8").Value = 0.1277189908446358
$ws.Range("P8").Value = 0.1277189908446358
$ws.Range("Q8").Value = 1.947159158262888
$ws.Range("R8").Value = 17.524432424366
$ws.Range("S8").Value = 0.06221025740983773
$ws.Range("T8").Value = 0.06221025740983773

# --- Row 9 (MuSCs -> FAPs) ---
$ws.Range("G9").Value = 0.4609646666666666
$ws.Range("I9").Value = 0.4870869789874365
$ws.Range("J9").Value = 0.4870869789874365
$ws.Range("O9").Value = 0.4492078640046304
$ws.Range("P9").Value = 0.4492078640046304
$ws.Range("Q9").Value = 6.848466313238666
$ws.Range("R9").Value = 61.63619681914799
$ws.Range("S9").Value = 0.2188033014154146
$ws.Range("T9").Value = 0.2188033014154146

# --- Row 10 (MuSCs -> MuSCs) ---
$ws.Range("G10").Value = 0.4609646666666666
$ws.Range("I10").Value = 0.4870869789874365
$ws.Range("J10").Value = 0.4870869789874365
$ws.Range("O10").Value = 0.4230731451507339
$ws.Range("P10").Value = 0.4230731451507338
$ws.Range("Q10").Value = 6.450025511064666
$ws.Range("S10").Value = 0.2060734201621842
$ws.Range("T10").Value = 0.2060734201621842
